$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 100.8373843333333
$ws.Range("H2").Value = 302.512153
$ws.Range("I2").Value = 0.6551985585448407
$ws.Range("J2").Value = 0.6551985585448408
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 78.370804
$ws.Range("N2").Value = 235.112412
$ws.Range("O2").Value = 0.9256182775132763
$ws.Range("P2").Value = 0.9256182775132761
$ws.Range("Q2").Value = 7902.706883460338
$ws.Range("R2").Value = 71124.36195114304
$ws.Range("S2").Value = 0.6064637611894569
$ws.Range("T2").Value = 0.6064637611894569
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 100.8373843333333
$ws.Range("H3").Value = 302.512153
$ws.Range("I3").Value = 0.6551985585448407
$ws.Range("J3").Value = 0.6551985585448408
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.309992333333333
$ws.Range("N3").Value = 6.929977
$ws.Range("O3").Value = 0.02728275091638557
$ws.Range("P3").Value = 0.02728275091638557
$ws.Range("Q3").Value = 232.9335847233868
$ws.Range("R3").Value = 2096.402262510481
$ws.Range("S3").Value = 0.01787561907355375
$ws.Range("T3").Value = 0.01787561907355376
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 100.8373843333333
$ws.Range("H4").Value = 302.512153
$ws.Range("I4").Value = 0.6551985585448407
$ws.Range("J4").Value = 0.6551985585448408
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.987804
$ws.Range("N4").Value = 11.963412
$ws.Range("O4").Value = 0.04709897157033827
$ws.Range("P4").Value = 0.04709897157033827
$ws.Range("Q4").Value = 402.119724594004
$ws.Range("R4").Value = 3619.077521346036
$ws.Range("S4").Value = 0.03085917828183007
$ws.Range("T4").Value = 0.03085917828183007
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.90796933333333
$ws.Range("H5").Value = 71.723908
$ws.Range("I5").Value = 0.1553438454249564
$ws.Range("J5").Value = 0.1553438454249564
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 78.370804
$ws.Range("N5").Value = 235.112412
$ws.Range("O5").Value = 0.9256182775132763
$ws.Range("P5").Value = 0.9256182775132761
$ws.Range("Q5").Value = 1873.686778660677
$ws.Range("R5").Value = 16863.1810079461
$ws.Range("S5").Value = 0.1437891026245368
$ws.Range("T5").Value = 0.1437891026245368
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.90796933333333
$ws.Range("H6").Value = 71.723908
$ws.Range("I6").Value = 0.1553438454249564
$ws.Range("J6").Value = 0.1553438454249564
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.309992333333333
$ws.Range("N6").Value = 6.929977
$ws.Range("O6").Value = 0.02728275091638557
$ws.Range("P6").Value = 0.02728275091638557
$ws.Range("Q6").Value = 55.22722586556844
$ws.Range("R6").Value = 497.045032790116
$ws.Range("S6").Value = 0.004238207441122587
$ws.Range("T6").Value = 0.004238207441122587
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.90796933333333
$ws.Range("H7").Value = 71.723908
$ws.Range("I7").Value = 0.1553438454249564
$ws.Range("J7").Value = 0.1553438454249564
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.987804
$ws.Range("N7").Value = 11.963412
$ws.Range("O7").Value = 0.04709897157033827
$ws.Range("P7").Value = 0.04709897157033827
$ws.Range("Q7").Value = 95.34029573934399
$ws.Range("R7").Value = 858.0626616540959
$ws.Range("S7").Value = 0.007316535359297045
$ws.Range("T7").Value = 0.007316535359297044
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 29.15819666666667
$ws.Range("H8").Value = 87.47459
$ws.Range("I8").Value = 0.1894575960302029
$ws.Range("J8").Value = 0.1894575960302029
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 78.370804
$ws.Range("N8").Value = 235.112412
$ws.Range("O8").Value = 0.9256182775132763
$ws.Range("P8").Value = 0.9256182775132761
$ws.Range("Q8").Value = 2285.151315956787
$ws.Range("R8").Value = 20566.36184361108
$ws.Range("S8").Value = 0.1753654136992826
$ws.Range("T8").Value = 0.1753654136992825
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 29.15819666666667
$ws.Range("H9").Value = 87.47459
$ws.Range("I9").Value = 0.1894575960302029
$ws.Range("J9").Value = 0.1894575960302029
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.309992333333333
$ws.Range("N9").Value = 6.929977
$ws.Range("O9").Value = 0.02728275091638557
$ws.Range("P9").Value = 0.02728275091638557
$ws.Range("Q9").Value = 67.35521075382556
$ws.Range("R9").Value = 606.19689678443
$ws.Range("S9").Value = 0.005168924401709225
$ws.Range("T9").Value = 0.005168924401709225
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 29.15819666666667
$ws.Range("H10").Value = 87.47459
$ws.Range("I10").Value = 0.1894575960302029
$ws.Range("J10").Value = 0.1894575960302029
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.987804
$ws.Range("N10").Value = 11.963412
$ws.Range("O10").Value = 0.04709897157033827
$ws.Range("P10").Value = 0.04709897157033827
$ws.Range("Q10").Value = 116.27717330012
$ws.Range("R10").Value = 1046.49455970108
$ws.Range("S10").Value = 0.008923257929211162
$ws.Range("T10").Value = 0.00892325792921116